# The sheet lists quarterly rows (A/B/C/D) for each year, in columns A:E
# (A: label, B: 产销率_累计值, C: 产销率比上年同期增减, D: 期末库存比年初增减,
#  E: 销售量_累计值). The "B" and "C" quarter rows within each year block were
# swapped (their A:E content exchanged), and columns F ("产销率") and G
# ("销售量") - which duplicated/derived data already present - were removed
# entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "C", "D", "E")

for ($base = 3; $base -le 63; $base += 4) {
    $rowB = $base
    $rowC = $base + 1

    foreach ($col in $cols) {
        $cellB = $ws.Range("$col$rowB")
        $cellC = $ws.Range("$col$rowC")
        $valB = $cellB.Value2
        $valC = $cellC.Value2
        $cellB.Value = $valC
        $cellC.Value = $valB
    }
}

# Remove the now-redundant F (产销率) and G (销售量) columns entirely.
$ws.Range("F1:G65").EntireColumn.Delete()

Write-Host "swap+delete complete"
